# Increase the gap between events dates, to allow new events to be
# inserted between the existing ones.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# confirmations
# ---------------------------------------------------------------------
$wsConfirmations = $wb.Worksheets.Item("confirmations")
$wsConfirmations.Range("A3").Value = 44576
$wsConfirmations.Range("A4").Value = 44592
$wsConfirmations.Range("A5").Value = 44593
$wsConfirmations.Range("A6").Value = 44607
$wsConfirmations.Range("A7").Value = 44620
$wsConfirmations.Range("A8").Value = 44620
$wsConfirmations.Range("A9").Value = 44621
$wsConfirmations.Range("A10").Value = 44635

# ---------------------------------------------------------------------
# trades
# ---------------------------------------------------------------------
$wsTrades = $wb.Worksheets.Item("trades")
$wsTrades.Range("A3").Value = 44576
$wsTrades.Range("A4").Value = 44576
$wsTrades.Range("A5").Value = 44592
$wsTrades.Range("A6").Value = 44592
$wsTrades.Range("A7").Value = 44593
$wsTrades.Range("A8").Value = 44593
$wsTrades.Range("A9").Value = 44607
$wsTrades.Range("A10").Value = 44620
$wsTrades.Range("A11").Value = 44620
$wsTrades.Range("A12").Value = 44620
$wsTrades.Range("A13").Value = 44621
$wsTrades.Range("A14").Value = 44635
[void]$wsTrades.Range("A2:A14").Select()

# ---------------------------------------------------------------------
# subscriptions
# ---------------------------------------------------------------------
$wsSubscriptions = $wb.Worksheets.Item("subscriptions")
$wsSubscriptions.Range("F2").Value = 44571
$wsSubscriptions.Range("G2").Value = 44572
$wsSubscriptions.Range("K2").Value = 44576
$wsSubscriptions.Range("A3").Value = 44652
$wsSubscriptions.Range("E3").Value = 44652
$wsSubscriptions.Range("F3").Value = 44661
$wsSubscriptions.Range("G3").Value = 44662

# ---------------------------------------------------------------------
# splits
# ---------------------------------------------------------------------
$wsSplits = $wb.Worksheets.Item("splits")
$wsSplits.Range("A2").Value = 44607
$wsSplits.Range("A3").Value = 44651
[void]$wsSplits.Range("A3").Select()

# ---------------------------------------------------------------------
# mergers
# ---------------------------------------------------------------------
$wsMergers = $wb.Worksheets.Item("mergers")
$wsMergers.Range("A2").Value = 44651
[void]$wsMergers.Range("A2").Select()

# ---------------------------------------------------------------------
# spinoffs
# ---------------------------------------------------------------------
$wsSpinoffs = $wb.Worksheets.Item("spinoffs")
$wsSpinoffs.Range("A2").Value = 44652

# ---------------------------------------------------------------------
# Make "confirmations" the active sheet/tab again.
# ---------------------------------------------------------------------
[void]$wsConfirmations.Select()
